$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update financial figures (Yearly financials refresh)
$ws.Range("D43").Value = 48200
$ws.Range("D45").Value = 18800
$ws.Range("D48").Value = 5500
$ws.Range("D49").Value = 45300

$ws.Range("D58").Value = 8800
$ws.Range("D59").Value = 56100
$ws.Range("D60").Value = 45400
$ws.Range("D61").Value = 5000

$ws.Range("D91").Value = -1200
$ws.Range("E91").Value = -2100
$ws.Range("F91").Value = -2900
$ws.Range("G91").Value = -800
$ws.Range("H91").Value = -1000
$ws.Range("I91").Value = -700
